# Update the "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets, which share the same rows of event data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 176
    4  = 12432
    6  = 143
    7  = 29
    10 = 195
    11 = 450
    16 = 374
    17 = 3619
    18 = 95
    19 = 940
    21 = 123
    22 = 50
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
